$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43; this shifts existing rows 43:51 down to 44:52
$ws.Rows.Item(43).Insert()

# Populate the new row 43 with a new price record for Guayaba (weekly update)
$ws.Range("A43").Value = 1
$ws.Range("B43").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C43").Value = "Arica y Parinacota"
$ws.Range("D43").Value = 44841
$ws.Range("E43").Value = 15
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100108
$ws.Range("H43").Value = "Tropicales y subtropicales"
$ws.Range("I43").Value = 100108001
$ws.Range("J43").Value = "Guayaba"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 1400
$ws.Range("O43").Value = 1500
$ws.Range("P43").Value = 1450
$ws.Range("Q43").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R43").Value = "Región de Arica y Parinacota"
$ws.Range("S43").Value = 1450
$ws.Range("T43").Value = 1
